$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..68 (date, colB, colC, colD)
$data = @(
    @("2018-10",108.9,100,100),
    @("2018-11",104.3,100.1,100.1),
    @("2018-12",101.7,100.3,100.3),
    @("2018-01",106.4502,99.2033,99.2033),
    @("2018-02",104.8,98.2,98.2),
    @("2018-03",105.1,95.7,95.7),
    @("2018-04",106.4,96.8,96.8),
    @("2018-05",108.2,98.5,98.5),
    @("2018-06",108.1,100.4,100.4),
    @("2018-07",106.7,100.4,100.4),
    @("2018-08",108.7,98.09999999999999,98.09999999999999),
    @("2018-09",111,98.7,98.7),
    @("2019-10",86.09999999999999,88,88),
    @("2019-11",87.3,88.3,88.3),
    @("2019-12",89,88.90000000000001,88.90000000000001),
    @("2019-01",99.8,98.5,98.5),
    @("2019-02",99.8,97.90000000000001,97.90000000000001),
    @("2019-03",99.5,97.59999999999999,97.59999999999999),
    @("2019-04",100.8,96,96),
    @("2019-05",99,95.90000000000001,95.90000000000001),
    @("2019-06",95.5,93.59999999999999,93.59999999999999),
    @("2019-07",95.7,92.59999999999999,92.59999999999999),
    @("2019-08",90.8,93.09999999999999,93.09999999999999),
    @("2019-09",86.40000000000001,90.3,90.3),
    @("2020-10",85.40000000000001,90.09999999999999,90.09999999999999),
    @("2020-11",88.5,94.40000000000001,94.40000000000001),
    @("2020-12",91.2,94.8,94.8),
    @("2020-01",90.59999999999999,88,88),
    @("2020-02",89.7,87.5,87.5),
    @("2020-03",87.5,88.40000000000001,88.40000000000001),
    @("2020-04",81,88.7,88.7),
    @("2020-05",81.3,88.40000000000001,88.40000000000001),
    @("2020-06",84.7,89.3,89.3),
    @("2020-07",82.8,88.8,88.8),
    @("2020-08",84.2,87.5,87.5),
    @("2020-09",83.59999999999999,88.59999999999999,88.59999999999999),
    @("2021-10",127.5,109.5,125.2),
    @("2021-11",125.8,107.1,120.7),
    @("2021-12",118.2,107.1,119.3),
    @("2021-01",94.2,91.7,103.2),
    @("2021-02",97.3,91.7,109.9),
    @("2021-03",109.4,92.2,115),
    @("2021-04",116.9,95.59999999999999,120.8),
    @("2021-05",117.9,95.5,120.4),
    @("2021-06",116.1,99.09999999999999,120.2),
    @("2021-07",120.7,104.7,121.8),
    @("2021-08",123.8,108,125),
    @("2021-09",123.5,107.7,125.4),
    @("2022-10",97.40000000000001,95.90000000000001,104.3),
    @("2022-11",94.7,95.90000000000001,101.8),
    @("2022-12",97.2,95.90000000000001,99.5),
    @("2022-01",114,107.1,114.6),
    @("2022-02",112.1,107.1,109.2),
    @("2022-03",104.4,100.6,105),
    @("2022-04",104.4,100.6,102),
    @("2022-05",106.1,100.6,104.4),
    @("2022-06",109.2,100.6,107.8),
    @("2022-07",104.7,100.6,108.2),
    @("2022-08",100.9,98.7,106.9),
    @("2022-09",102.3,97.2,106),
    @("2023-01",97.40000000000001,95.90000000000001,98.59999999999999),
    @("2023-02",97.3,93.09999999999999,100.4),
    @("2023-03",95.5,99.09999999999999,99.5),
    @("2023-04",96.40000000000001,99.09999999999999,99.59999999999999),
    @("2023-05",95,99.09999999999999,97.90000000000001),
    @("2023-06",92.7,99.09999999999999,95.59999999999999),
    @("2023-07",93.8,99.09999999999999,96.40000000000001)
)

$existingLastRow = 49
$targetLastRow = 1 + $data.Count

# Copy formatting (style) from row 2 (A2:D2) into any brand-new rows
# that do not yet exist in the original sheet (rows 50..targetLastRow).
if ($targetLastRow -gt $existingLastRow) {
    $ws.Range("A2:D2").Copy()
    for ($r = $existingLastRow + 1; $r -le $targetLastRow; $r++) {
        $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)
    }
    $excel.CutCopyMode = 0
}

# Write values row by row
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "done, rows written:" $data.Count
